$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "sheet"

$ws.Range("A2").Value = "BR"
$ws.Range("D2").Value = 20.758449090758884

$ws.Range("A3").Value = "NE"
$ws.Range("D3").Value = 26.071652978380218

$ws.Range("A4").Value = "SE"
$ws.Range("D4").Value = 18.220477881476292
$ws.Range("E4").Value = 19

$ws.Range("A5").Value = "AM"
$ws.Range("D5").Value = 28.529403728631156
$ws.Range("E5").Value = 6

$ws.Range("A6").Value = "BA"
$ws.Range("D6").Value = 30.627262240958647
$ws.Range("E6").Value = 5

$ws.Range("A7").Value = "CE"
$ws.Range("D7").Value = 30.964663250456407
$ws.Range("E7").Value = 4

$ws.Range("A8").Value = "AL"
$ws.Range("D8").Value = 32.289416974030253
$ws.Range("E8").Value = 3

$ws.Range("A9").Value = "AP"
$ws.Range("D9").Value = 33.344741710310032
$ws.Range("E9").Value = 2

$ws.Range("A10").Value = "PE"
$ws.Range("D10").Value = 34.046323845117534
$ws.Range("E10").Value = 1

$ws.Range("A2:E10").Select()
$excel.ActiveCell = $ws.Range("A2")
